$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Locate the "console.log(...)" paragraph (it's the 3rd paragraph: an empty
# paragraph, the "// Aqui puedes..." comment paragraph, then this one).
# --------------------------------------------------------------------------
$codePara = $d.Paragraphs.Item(3)

# The paragraph currently looks like:
#   [gramStart]<r>console.log(</r>[gramEnd]<r>"Bienvenido a la Inmobiliaria </r>[spellStart]<r>Marin</r>[spellEnd]<r>");</r>
# and must become:
#   <r>console.log("Bienvenido a la Inmobiliaria Iovanni </r>[spellStart]<r>Marin</r>[spellEnd]<r>");</r>
#
# i.e. the leading grammar-check markers around "console.log(" must be
# dropped and the first two runs merged into a single run with the new
# wording, while the spell-check markers around "Marin" (and the trailing
# run) must stay untouched.
#
# A plain Find/Replace across that span leaves the *leading* proofErr mark
# behind because it sits exactly on the paragraph/range boundary (marks
# strictly *inside* a replaced range are dropped, but one sitting right at
# the edge survives). To make the leading mark an "interior" position, we
# first insert a one-character placeholder immediately before it, shifting
# it off the boundary, then do the real replacement across the placeholder
# *and* the old text. Using plain Range.Text assignment (not Find.Execute)
# also avoids AutoCorrect turning the straight quote into a curly one.

$codePara.Range.InsertBefore("X")

$codePara = $d.Paragraphs.Item(3)
$oldPrefix = 'console.log("Bienvenido a la Inmobiliaria '
$newPrefix = 'console.log("Bienvenido a la Inmobiliaria Iovanni '
$start = $codePara.Range.Start
$prefixRange = $d.Range($start, $start + 1 + $oldPrefix.Length)
$prefixRange.Text = $newPrefix

# --------------------------------------------------------------------------
# Remove the empty trailing paragraph that sits right before the sectPr.
# Deleting the final paragraph mark on its own is a no-op (it's the body's
# sentinel mark), so instead we delete the combined range spanning the code
# paragraph's own mark *and* the empty paragraph's mark. That collapses the
# empty paragraph away while letting the code paragraph keep its own
# identity/properties.
# --------------------------------------------------------------------------
$codePara = $d.Paragraphs.Item(3)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$mergeRange = $d.Range($codePara.Range.End - 1, $lastPara.Range.End)
$mergeRange.Delete()
